$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("MegaMenuInfo")
$ws.Range("A3").Value = "/espanol/investigacion"
